# Update "F" column (想去人数 / interest count) values on the 展览 and
# 全部类型 sheets to match the freshly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 636
$ws1.Range("F3").Value  = 2214
$ws1.Range("F5").Value  = 13243
$ws1.Range("F7").Value  = 117
$ws1.Range("F8").Value  = 518
$ws1.Range("F12").Value = 13793
$ws1.Range("F13").Value = 14425
$ws1.Range("F22").Value = 1097
$ws1.Range("F25").Value = 5475
$ws1.Range("F26").Value = 941
$ws1.Range("F27").Value = 430
$ws1.Range("F28").Value = 339
$ws1.Range("F30").Value = 84

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 636
$ws4.Range("F3").Value  = 2214
$ws4.Range("F5").Value  = 13243
$ws4.Range("F8").Value  = 117
$ws4.Range("F9").Value  = 518
$ws4.Range("F13").Value = 13793
$ws4.Range("F14").Value = 14425
$ws4.Range("F23").Value = 1097
$ws4.Range("F26").Value = 5475
$ws4.Range("F28").Value = 430
$ws4.Range("F29").Value = 339
$ws4.Range("F31").Value = 84
